$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. Insert a new row at
# position 393 (shifting the existing rows 393-477 down to 394-478) and
# populate it with the new record's data.
$ws.Rows(393).Insert()

$ws.Range("A393").Value = 10
$ws.Range("B393").Value = "Vega Modelo de Temuco"
$ws.Range("C393").Value = "La Araucanía"
$ws.Range("D393").Value = 44159
$ws.Range("E393").Value = 9
$ws.Range("F393").Value = 100114001
$ws.Range("G393").Value = "Papa"
$ws.Range("H393").Value = "Pehuenche"
$ws.Range("I393").Value = "1a nueva(o)"
$ws.Range("J393").Value = 200
$ws.Range("K393").Value = 14000
$ws.Range("L393").Value = 14000
$ws.Range("M393").Value = 14000
$ws.Range("N393").Value = "$/saco 25 kilos"
$ws.Range("O393").Value = "Provincia de Cautín"
$ws.Range("P393").Value = 560
$ws.Range("Q393").Value = 25
$ws.Range("R393").Value = "Hortaliza"
